$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column is added as column H, mirroring the header style used by
# the other header cells (e.g. G1: bold font, border, centered/top aligned).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
